$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Pros1", "Tyro3", "ECs", 3, 1, 15.33750333333333, 46.01251, 0.1440483515229198, 0.1440483515229198, 1, 0.3333333333333333, 0.096929, 0.290787, 0.01940247578118474, 0.01940247578118474, 1.486648860596667, 13.37983974537, 0.002794894651743038, 0.002794894651743038),
  @("ECs", "Pros1", "Tyro3", "FAPs", 3, 1, 15.33750333333333, 46.01251, 0.1440483515229198, 0.1440483515229198, 3, 1, 3.954952666666667, 11.864858, 0.7916709481242146, 0.7916709481242146, 60.65909970817555, 545.93189737358, 0.1140388950258801, 0.1140388950258801),
  @("ECs", "Pros1", "Tyro3", "sCs", 3, 1, 15.33750333333333, 46.01251, 0.1440483515229198, 0.1440483515229198, 3, 1, 0.9438210000000001, 2.831463, 0.1889265760946008, 0.1889265760946008, 14.47585773357, 130.28271960213, 0.02721456184529672, 0.02721456184529672),
  @("FAPs", "Pros1", "Tyro3", "ECs", 3, 1, 64.92210766666668, 194.766323, 0.609742171429695, 0.6097421714296949, 1, 0.3333333333333333, 0.096929, 0.290787, 0.01940247578118474, 0.01940247578118474, 6.292834974022334, 56.63551476620101, 0.01183050771393165, 0.01183050771393165),
  @("FAPs", "Pros1", "Tyro3", "FAPs", 3, 1, 64.92210766666668, 194.766323, 0.609742171429695, 0.6097421714296949, 3, 1, 3.954952666666667, 11.864858, 0.7916709481242146, 0.7916709481242146, 256.7638628419038, 2310.874765577134, 0.482715162967064, 0.4827151629670639),
  @("FAPs", "Pros1", "Tyro3", "sCs", 3, 1, 64.92210766666668, 194.766323, 0.609742171429695, 0.6097421714296949, 3, 1, 0.9438210000000001, 2.831463, 0.1889265760946008, 0.1889265760946008, 61.27484858006102, 551.4736372205491, 0.1151965007486994, 0.1151965007486994),
  @("M2", "Pros1", "Tyro3", "ECs", 3, 1, 10.67805633333333, 32.034169, 0.1002872748488753, 0.1002872748488753, 1, 0.3333333333333333, 0.096929, 0.290787, 0.01940247578118474, 0.01940247578118474, 1.035013322333667, 9.315119901002999, 0.00194582142141632, 0.00194582142141632),
  @("M2", "Pros1", "Tyro3", "FAPs", 3, 1, 10.67805633333333, 32.034169, 0.1002872748488753, 0.1002872748488753, 3, 1, 3.954952666666667, 11.864858, 0.7916709481242146, 0.7916709481242146, 42.23120737033355, 380.080866333002, 0.07939452196440278, 0.07939452196440278),
  @("M2", "Pros1", "Tyro3", "sCs", 3, 1, 10.67805633333333, 32.034169, 0.1002872748488753, 0.1002872748488753, 3, 1, 0.9438210000000001, 2.831463, 0.1889265760946008, 0.1889265760946008, 10.078173806583, 90.703564259247, 0.01894693146305618, 0.01894693146305618),
  @("sCs", "Pros1", "Tyro3", "ECs", 3, 1, 15.537021, 46.611063, 0.14592220219851, 0.14592220219851, 1, 0.3333333333333333, 0.096929, 0.290787, 0.01940247578118474, 0.01940247578118474, 1.505987908509, 13.553891176581, 0.002831251994093733, 0.002831251994093732),
  @("sCs", "Pros1", "Tyro3", "FAPs", 3, 1, 15.537021, 46.611063, 0.14592220219851, 0.14592220219851, 3, 1, 3.954952666666667, 11.864858, 0.7916709481242146, 0.7916709481242146, 61.448182636006, 553.0336437240541, 0.1155223681668678, 0.1155223681668677),
  @("sCs", "Pros1", "Tyro3", "sCs", 3, 1, 15.537021, 46.611063, 0.14592220219851, 0.14592220219851, 3, 1, 0.9438210000000001, 2.831463, 0.1889265760946008, 0.1889265760946008, 14.664166697241, 131.977500275169, 0.02756858203754853, 0.02756858203754852)
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
  }
}
